# RPA-93: Lage intitieringsfil og rutine
#
# Changes applied to the "Vasklister" header sheet (Ark1):
#   1. A new "Kommune_Nr" column is inserted at column C (header row 1),
#      pushing FNR (previously C) and everything after it one column to
#      the right.
#   2. The old "_1_inntekt" column (which ends up at column E after the
#      insert above) is removed, shifting everything after it back left
#      by one column - so the sheet still ends up with 20 columns (A:T).
#   3. The sample data in row 2 is refreshed with new FNR / Sak_Nr values
#      (and the mirrored value in C2 is kept in sync with the new FNR).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header row, column C onward (A and B - debitor_ident / Sak_Nr - stay put).
$headers = @(
    "Kommune_Nr",
    "FNR",
    "_1_isProcentTrekk",
    "_1_tbNrEllerPct",
    "_2_numVoksne",
    "_23_isDelUtgift",
    "_2_barn05",
    "_2_barn610",
    "_2_barn1118",
    "_3_region",
    "_4_kontantstøtte",
    "_4_barnetrygd",
    "_4_andreUtgift",
    "_4_orgNr",
    "_4_firmaNavn",
    "_4_firmaPostAddresse",
    "_4_firmaPostNr",
    "_5_ErrorCode"
)

$col = 3
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# Refresh the sample row underneath the headers.
$ws.Range("A2").Value = "13098245418"
$ws.Range("B2").Value = "313054"
$ws.Range("C2").Value = "13098245418"
